# Auto-generated script to apply commit changes to Cactuar Profits workbook
$wb = $excel.ActiveWorkbook

# --- Sheet: ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H106").Value = 22224422
$ws.Range("I106").Value = 25642642
$ws.Range("J106").Value = 5990
$ws.Range("K106").Value = 25642642
$ws.Range("L106").Value = 5990
$ws.Range("M106").Value = -25642011
$ws.Range("N106").Value = -7252
$ws.Range("H113").Value = 3486.182
$ws.Range("I113").Value = 2713.7144
$ws.Range("K113").Value = 2713.7144
$ws.Range("M113").Value = 540.2856000000002
$ws.Range("H121").Value = 4188.143
$ws.Range("J121").Value = 4188.143
$ws.Range("L121").Value = 12564.429
$ws.Range("N121").Value = -16058.429
$ws.Range("H132").Value = 18957.4
$ws.Range("I132").Value = 9236.833000000001
$ws.Range("J132").Value = 23686.324
$ws.Range("K132").Value = 27710.499
$ws.Range("L132").Value = 71058.97200000001
$ws.Range("M132").Value = -25180.499
$ws.Range("N132").Value = -76118.97200000001
$ws.Range("H141").Value = 6398.3335
$ws.Range("I141").Value = 6398.3335
$ws.Range("K141").Value = 19195.0005
$ws.Range("M141").Value = -14015.0005

# --- Sheet: ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 15588.705
$ws.Range("I32").Value = 15987.06
$ws.Range("K32").Value = 15987.06
$ws.Range("M32").Value = -15700.06
$ws.Range("H74").Value = 9616332
$ws.Range("I74").Value = 14706678
$ws.Range("K74").Value = 14706678
$ws.Range("M74").Value = -14705804
$ws.Range("H77").Value = 9616332
$ws.Range("I77").Value = 14706678
$ws.Range("K77").Value = 73533390
$ws.Range("M77").Value = -73529022
$ws.Range("H97").Value = 551.5833
$ws.Range("J97").Value = 637.5
$ws.Range("L97").Value = 637.5
$ws.Range("N97").Value = -1629.5
$ws.Range("H110").Value = 1364041.5
$ws.Range("I110").Value = 2916767.5
$ws.Range("K110").Value = 2916767.5
$ws.Range("M110").Value = -2914722.5
$ws.Range("H122").Value = 3215.1562
$ws.Range("I122").Value = 1962.2916
$ws.Range("K122").Value = 5886.8748
$ws.Range("M122").Value = -3436.8748
$ws.Range("H132").Value = 14566.413
$ws.Range("I132").Value = 17792.854
$ws.Range("J132").Value = 5424.8335
$ws.Range("K132").Value = 53378.562
$ws.Range("L132").Value = 16274.5005
$ws.Range("M132").Value = -50848.562
$ws.Range("N132").Value = -21334.5005

# --- Sheet: BSM ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H20").Value = 2663.4285
$ws.Range("I20").Value = 3036.5386
$ws.Range("J20").Value = 2057.125
$ws.Range("K20").Value = 3036.5386
$ws.Range("L20").Value = 2057.125
$ws.Range("M20").Value = -2789.5386
$ws.Range("N20").Value = -2551.125
$ws.Range("H105").Value = 100001800
$ws.Range("I105").Value = 125001510
$ws.Range("J105").Value = 2944
$ws.Range("K105").Value = 125001510
$ws.Range("L105").Value = 2944
$ws.Range("M105").Value = -124999763
$ws.Range("N105").Value = -6438
$ws.Range("H134").Value = 2894.913
$ws.Range("I134").Value = 1107.4375
$ws.Range("J134").Value = 6980.5713
$ws.Range("K134").Value = 3322.3125
$ws.Range("L134").Value = 20941.7139
$ws.Range("M134").Value = -787.3125
$ws.Range("N134").Value = -26011.7139

# --- Sheet: CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H99").Value = 18909.637
$ws.Range("J99").Value = 9999.556
$ws.Range("L99").Value = 9999.556
$ws.Range("N99").Value = -12995.556
$ws.Range("H126").Value = 18909.637
$ws.Range("J126").Value = 9999.556
$ws.Range("L126").Value = 29998.668
$ws.Range("N126").Value = -34938.66800000001
$ws.Range("H132").Value = 11503282
$ws.Range("I132").Value = 13899324
$ws.Range("K132").Value = 41697972
$ws.Range("M132").Value = -41695442
$ws.Range("H141").Value = 104761.75
$ws.Range("J141").Value = 110249.27
$ws.Range("L141").Value = 110249.27
$ws.Range("N141").Value = -120609.27

# --- Sheet: CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H4").Value = 5329646.5
$ws.Range("I4").Value = 1921674.4
$ws.Range("J4").Value = 17257550
$ws.Range("K4").Value = 5765023.199999999
$ws.Range("L4").Value = 51772650
$ws.Range("M4").Value = -5764911.199999999
$ws.Range("N4").Value = -51772874
$ws.Range("H113").Value = 729.7895
$ws.Range("I113").Value = 704.8570999999999
$ws.Range("J113").Value = 744.3333
$ws.Range("K113").Value = 2114.5713
$ws.Range("L113").Value = 2232.9999
$ws.Range("M113").Value = 55.42870000000039
$ws.Range("N113").Value = -6572.9999
$ws.Range("H134").Value = 12566.385
$ws.Range("I134").Value = 3713.111
$ws.Range("K134").Value = 11139.333
$ws.Range("M134").Value = -6069.332999999999

# --- Sheet: GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H63").Value = 82332.664
$ws.Range("I63").Value = 0
$ws.Range("J63").Value = 82332.664
$ws.Range("K63").Value = 0
$ws.Range("L63").Value = 82332.664
$ws.Range("M63").ClearContents()
$ws.Range("N63").Value = -83704.664
$ws.Range("H66").Value = 82332.664
$ws.Range("I66").Value = 0
$ws.Range("J66").Value = 82332.664
$ws.Range("K66").Value = 0
$ws.Range("L66").Value = 246997.992
$ws.Range("M66").ClearContents()
$ws.Range("N66").Value = -253861.992
$ws.Range("H80").Value = 1329165.4
$ws.Range("I80").Value = 3127999.8
$ws.Range("J80").Value = 44283.715
$ws.Range("K80").Value = 3127999.8
$ws.Range("L80").Value = 44283.715
$ws.Range("M80").Value = -3127001.8
$ws.Range("N80").Value = -46279.715
$ws.Range("H83").Value = 1329165.4
$ws.Range("I83").Value = 3127999.8
$ws.Range("J83").Value = 44283.715
$ws.Range("K83").Value = 15639999
$ws.Range("L83").Value = 221418.575
$ws.Range("M83").Value = -15635007
$ws.Range("N83").Value = -231402.575
$ws.Range("H97").Value = 550.62964
$ws.Range("I97").Value = 522.4
$ws.Range("J97").Value = 631.2857
$ws.Range("K97").Value = 522.4
$ws.Range("L97").Value = 631.2857
$ws.Range("M97").Value = -26.39999999999998
$ws.Range("N97").Value = -1623.2857
$ws.Range("H107").Value = 7937341.5
$ws.Range("J107").Value = 999.6667
$ws.Range("L107").Value = 999.6667
$ws.Range("N107").Value = -4839.6667
$ws.Range("H122").Value = 327289.56
$ws.Range("I122").Value = 460639.53
$ws.Range("K122").Value = 1381918.59
$ws.Range("M122").Value = -1379468.59
$ws.Range("H132").Value = 105991.7
$ws.Range("I132").Value = 157683.77
$ws.Range("J132").Value = 9992.143
$ws.Range("K132").Value = 473051.3099999999
$ws.Range("L132").Value = 29976.429
$ws.Range("M132").Value = -470521.3099999999
$ws.Range("N132").Value = -35036.429
$ws.Range("H138").Value = 100429
$ws.Range("I138").Value = 0
$ws.Range("J138").Value = 100429
$ws.Range("K138").Value = 0
$ws.Range("L138").Value = 100429
$ws.Range("M138").ClearContents()
$ws.Range("N138").Value = -110709
$ws.Range("H140").Value = 29593
$ws.Range("J140").Value = 29593
$ws.Range("L140").Value = 29593
$ws.Range("N140").Value = -39953

# --- Sheet: LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 4309.8184
$ws.Range("I7").Value = 2095.2856
$ws.Range("K7").Value = 2095.2856
$ws.Range("M7").Value = -1983.2856
$ws.Range("H40").Value = 13891480
$ws.Range("I40").Value = 2256.9285
$ws.Range("K40").Value = 2256.9285
$ws.Range("M40").Value = -2120.9285
$ws.Range("H68").Value = 2527050.5
$ws.Range("J68").Value = 2666
$ws.Range("L68").Value = 2666
$ws.Range("N68").Value = -4164
$ws.Range("H71").Value = 2527050.5
$ws.Range("J71").Value = 2666
$ws.Range("L71").Value = 13330
$ws.Range("N71").Value = -20818
$ws.Range("H122").Value = 57147652
$ws.Range("I122").Value = 100003210
$ws.Range("K122").Value = 300009630
$ws.Range("M122").Value = -300007180
$ws.Range("H126").Value = 4309.8184
$ws.Range("I126").Value = 2095.2856
$ws.Range("K126").Value = 6285.8568
$ws.Range("M126").Value = -3815.8568
$ws.Range("H132").Value = 4558.7144
$ws.Range("I132").Value = 3485.818
$ws.Range("J132").Value = 6374.385
$ws.Range("K132").Value = 10457.454
$ws.Range("L132").Value = 19123.155
$ws.Range("M132").Value = -7927.454000000002
$ws.Range("N132").Value = -24183.155
$ws.Range("H136").Value = 4298.2
$ws.Range("I136").Value = 3007.9
$ws.Range("K136").Value = 9023.700000000001
$ws.Range("M136").Value = -6473.700000000001
$ws.Range("H139").Value = 79357.5
$ws.Range("J139").Value = 79357.5
$ws.Range("L139").Value = 79357.5
$ws.Range("N139").Value = -89637.5

# --- Sheet: WVR ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H2").Value = 3014.9092
$ws.Range("I2").Value = 2816.4
$ws.Range("K2").Value = 2816.4
$ws.Range("M2").Value = -2704.4
$ws.Range("H62").Value = 7526
$ws.Range("I62").Value = 4407.5
$ws.Range("K62").Value = 4407.5
$ws.Range("M62").Value = -3783.5
$ws.Range("H65").Value = 7526
$ws.Range("I65").Value = 4407.5
$ws.Range("K65").Value = 22037.5
$ws.Range("M65").Value = -18917.5
$ws.Range("H107").Value = 2466.3572
$ws.Range("I107").Value = 2556.5
$ws.Range("J107").Value = 2304.1
$ws.Range("K107").Value = 7669.5
$ws.Range("L107").Value = 6912.299999999999
$ws.Range("M107").Value = -5749.5
$ws.Range("N107").Value = -10752.3
$ws.Range("H126").Value = 2420.8948
$ws.Range("I126").Value = 2965.4285
$ws.Range("J126").Value = 2103.25
$ws.Range("K126").Value = 8896.2855
$ws.Range("L126").Value = 6309.75
$ws.Range("M126").Value = -6426.2855
$ws.Range("N126").Value = -11249.75
$ws.Range("H132").Value = 23150318
$ws.Range("I132").Value = 3472859.5
$ws.Range("J132").Value = 62505236
$ws.Range("K132").Value = 10418578.5
$ws.Range("L132").Value = 187515708
$ws.Range("M132").Value = -10416048.5
$ws.Range("N132").Value = -187520768
$ws.Range("H136").Value = 9500.794
$ws.Range("I136").Value = 2297.353
$ws.Range("J136").Value = 11901.941
$ws.Range("K136").Value = 6892.059
$ws.Range("L136").Value = 35705.823
$ws.Range("M136").Value = -4342.059
$ws.Range("N136").Value = -40805.823
